$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = 5.5
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.25
$ws.Range("R5").Value = 1.75
$ws.Range("G11").Value = 4.5
$ws.Range("M11").Value = 1.09
$ws.Range("N11").Value = 7
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 1.65
$ws.Range("W11").Value = 10.25
$ws.Range("X11").Value = 24
$ws.Range("AB11").Value = 60
$ws.Range("AD11").Value = 6.4
$ws.Range("AI11").Value = 8.5
$ws.Range("AU11").Value = 7.7
$ws.Range("AW11").Value = 3.5
$ws.Range("I12").Value = 4.05
$ws.Range("L12").Value = 4.25
$ws.Range("O12").Value = 1.2
$ws.Range("P12").Value = 3.6
$ws.Range("W12").Value = 8.25
$ws.Range("Y12").Value = 8.25
$ws.Range("AA12").Value = 12.5
$ws.Range("AB12").Value = 22
$ws.Range("AG12").Value = 14.5
$ws.Range("AP12").Value = 16
$ws.Range("AQ12").Value = 28
$ws.Range("AX12").Value = 21
$ws.Range("BB12").Value = 250
$ws.Range("G26").Value = 1.38
$ws.Range("I26").Value = 7
$ws.Range("S26").Value = 1.18
$ws.Range("T26").Value = 4.5
$ws.Range("BD26").Value = 176
$ws.Range("G29").Value = 2.77
$ws.Range("H29").Value = 3.55
$ws.Range("I29").Value = 2.2
$ws.Range("X29").Value = 15
$ws.Range("Y29").Value = 10.25
$ws.Range("AA29").Value = 22
$ws.Range("AD29").Value = 7
$ws.Range("AK29").Value = 17
$ws.Range("AM29").Value = 400
$ws.Range("AN29").Value = 4.85
$ws.Range("AP29").Value = 21
$ws.Range("AR29").Value = 90
$ws.Range("G30").Value = 1.78
$ws.Range("I30").Value = 3.6
$ws.Range("J30").Value = 2.27
$ws.Range("L30").Value = 3.85
$ws.Range("N30").Value = 9.5
$ws.Range("O30").Value = 1.15
$ws.Range("P30").Value = 4.65
$ws.Range("Q30").Value = 1.47
$ws.Range("R30").Value = 2.5
$ws.Range("T30").Value = 3.5
$ws.Range("V30").Value = 2.42
$ws.Range("X30").Value = 11
$ws.Range("Z30").Value = 16.5
$ws.Range("AA30").Value = 12.5
$ws.Range("AC30").Value = 9.5
$ws.Range("AE30").Value = 12
$ws.Range("AF30").Value = 37
$ws.Range("AG30").Value = 17
$ws.Range("AI30").Value = 12.5
$ws.Range("AK30").Value = 27
$ws.Range("AL30").Value = 26
$ws.Range("AN30").Value = 4.1
$ws.Range("AO30").Value = 8.75
$ws.Range("AP30").Value = 14
$ws.Range("AQ30").Value = 26
$ws.Range("AR30").Value = 45
$ws.Range("AT30").Value = 3.5
$ws.Range("AU30").Value = 6.5
$ws.Range("AW30").Value = 6.1
$ws.Range("AX30").Value = 18.5
$ws.Range("AY30").Value = 19.5
$ws.Range("AZ30").Value = 80
$ws.Range("BA30").Value = 90
$ws.Range("BB30").Value = 175
$ws.Range("G31").Value = 1.65
$ws.Range("H31").Value = 3.65
$ws.Range("I31").Value = 4.55
$ws.Range("J31").Value = 2.22
$ws.Range("L31").Value = 4.9
$ws.Range("P31").Value = 3.25
$ws.Range("Q31").Value = 1.88
$ws.Range("R31").Value = 1.85
$ws.Range("S31").Value = 1.39
$ws.Range("U31").Value = 1.87
$ws.Range("V31").Value = 1.85
$ws.Range("W31").Value = 6.7
$ws.Range("X31").Value = 7.6
$ws.Range("Z31").Value = 12.5
$ws.Range("AD31").Value = 7.2
$ws.Range("AE31").Value = 16.5
$ws.Range("AF31").Value = 80
$ws.Range("AH31").Value = 27
$ws.Range("AI31").Value = 15.5
$ws.Range("AJ31").Value = 80
$ws.Range("AL31").Value = 50
$ws.Range("AM31").Value = 700
$ws.Range("AN31").Value = 3.5
$ws.Range("AO31").Value = 8.25
$ws.Range("AP31").Value = 18
$ws.Range("AQ31").Value = 27
$ws.Range("AR31").Value = 60
$ws.Range("AU31").Value = 7.8
$ws.Range("AW31").Value = 6.3
$ws.Range("AX31").Value = 27
$ws.Range("BA31").Value = 200
$ws.Range("G32").Value = 2.18
$ws.Range("H32").Value = 3.2
$ws.Range("I32").Value = 3.05
$ws.Range("J32").Value = 2.77
$ws.Range("K32").Value = 2.1
$ws.Range("L32").Value = 3.65
$ws.Range("P32").Value = 3.35
$ws.Range("T32").Value = 2.75
$ws.Range("U32").Value = 1.65
$ws.Range("V32").Value = 2.1
$ws.Range("X32").Value = 11.5
$ws.Range("Y32").Value = 8.75
$ws.Range("Z32").Value = 22
$ws.Range("AA32").Value = 17
$ws.Range("AB32").Value = 25
$ws.Range("AH32").Value = 16.5
$ws.Range("AI32").Value = 10.75
$ws.Range("AJ32").Value = 40
$ws.Range("AK32").Value = 26
$ws.Range("AL32").Value = 32
$ws.Range("AN32").Value = 4.2
$ws.Range("AO32").Value = 11.25
$ws.Range("AP32").Value = 18.5
$ws.Range("AR32").Value = 70
$ws.Range("AS32").Value = 200
$ws.Range("AT32").Value = 2.75
$ws.Range("AU32").Value = 6.8
$ws.Range("AW32").Value = 5.1
$ws.Range("AX32").Value = 17.5
$ws.Range("AY32").Value = 24
$ws.Range("AZ32").Value = 80
$ws.Range("BA32").Value = 120
$ws.Range("BB32").Value = 300
$ws.Range("G33").Value = 2.9
$ws.Range("J33").Value = 3.45
$ws.Range("K33").Value = 2.22
$ws.Range("O33").Value = 1.26
$ws.Range("P33").Value = 3.5
$ws.Range("Q33").Value = 1.78
$ws.Range("R33").Value = 1.98
$ws.Range("W33").Value = 10.25
$ws.Range("Y33").Value = 10.75
$ws.Range("AA33").Value = 24
$ws.Range("AG33").Value = 8.5
$ws.Range("AV33").Value = 60
$ws.Range("AW33").Value = 4.15
$ws.Range("AY33").Value = 18.5
$ws.Range("G34").Value = 2.35
$ws.Range("H34").Value = 3.15
$ws.Range("I34").Value = 2.8
$ws.Range("J34").Value = 2.95
$ws.Range("K34").Value = 2.1
$ws.Range("L34").Value = 3.4
$ws.Range("T34").Value = 2.72
$ws.Range("U34").Value = 1.75
$ws.Range("V34").Value = 1.95
$ws.Range("W34").Value = 7.9
$ws.Range("X34").Value = 11.5
$ws.Range("Y34").Value = 9.25
$ws.Range("Z34").Value = 24
$ws.Range("AA34").Value = 19.5
$ws.Range("AB34").Value = 29
$ws.Range("AG34").Value = 8.75
$ws.Range("AH34").Value = 14.5
$ws.Range("AI34").Value = 10.25
$ws.Range("AJ34").Value = 35
$ws.Range("AK34").Value = 25
$ws.Range("AL34").Value = 35
$ws.Range("AN34").Value = 4.3
$ws.Range("AO34").Value = 12.5
$ws.Range("AQ34").Value = 50
$ws.Range("AR34").Value = 80
$ws.Range("AT34").Value = 2.72
$ws.Range("AW34").Value = 4.8
$ws.Range("AX34").Value = 15.5
$ws.Range("AY34").Value = 23
$ws.Range("AZ34").Value = 70
$ws.Range("BB34").Value = 300
$ws.Range("G37").Value = 2.35
$ws.Range("H37").Value = 2.9
$ws.Range("J37").Value = 2.9
$ws.Range("L37").Value = 3.7
$ws.Range("M37").Value = 1.02
$ws.Range("N37").Value = 7.1
$ws.Range("T37").Value = 2.47
$ws.Range("U37").Value = 1.75
$ws.Range("V37").Value = 1.85
$ws.Range("W37").Value = 7.2
$ws.Range("X37").Value = 11.5
$ws.Range("Y37").Value = 9
$ws.Range("AA37").Value = 20
$ws.Range("AB37").Value = 30
$ws.Range("AC37").Value = 7.6
$ws.Range("AD37").Value = 5.6
$ws.Range("AG37").Value = 8.25
$ws.Range("AH37").Value = 16
$ws.Range("AK37").Value = 30
$ws.Range("AL37").Value = 40
$ws.Range("AN37").Value = 4.2
$ws.Range("AO37").Value = 12
$ws.Range("AP37").Value = 19
$ws.Range("AQ37").Value = 50
$ws.Range("AR37").Value = 80
$ws.Range("AT37").Value = 2.45
$ws.Range("AU37").Value = 6.6
$ws.Range("AV37").Value = 55
$ws.Range("AX37").Value = 17.5
$ws.Range("AY37").Value = 24
$ws.Range("AZ37").Value = 90
$ws.Range("BA37").Value = 120
